{"js": "const pairs = [\n  [\"2024-04-20 Saturday\", \"2024-04-21 Sunday\"],\n  [\"34\u00d785=2890\", \"92\u00d752=4784\"],\n  [\"84\u00d731=2604\", \"11\u00d746=506\"],\n  [\"91\u00d730=2730\", \"64\u00d769=4416\"],\n  [\"98\u00d779=7742\", \"47\u00d726=1222\"],\n  [\"59\u00d780=4720\", \"54\u00d787=4698\"],\n  [\"34\u00d740=1360\", \"15\u00d714=210\"],\n  [\"92\u00d716=1472\", \"48\u00d712=576\"],\n  [\"88\u00d795=8360\", \"56\u00d783=4648\"],\n  [\"60\u00d783=4980\", \"66\u00d719=1254\"],\n  [\"67\u00d784=5628\", \"90\u00d713=1170\"],\n  [\"64\u00d752=3328\", \"92\u00d725=2300\"],\n  [\"53\u00d729=1537\", \"95\u00d789=8455\"],\n  [\"38\u00d751=1938\", \"19\u00d735=665\"],\n  [\"26\u00d794=2444\", \"25\u00d784=2100\"],\n  [\"59\u00d796=5664\", \"82\u00d723=1886\"],\n  [\"42\u00d776=3192\", \"46\u00d749=2254\"],\n  [\"45\u00d725=1125\", \"84\u00d777=6468\"],\n  [\"38\u00d793=3534\", \"89\u00d733=2937\"],\n  [\"62\u00d785=5270\", \"60\u00d739=2340\"],\n  [\"47\u00d753=2491\", \"17\u00d756=952\"],\n  [\"51\u00d767=3417\", \"62\u00d790=5580\"],\n  [\"61\u00d749=2989\", \"95\u00d728=2660\"],\n  [\"29\u00d745=1305\", \"36\u00d737=1332\"],\n  [\"19\u00d746=874\", \"35\u00d798=3430\"],\n  [\"59\u00d787=5133\", \"27\u00d776=2052\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Search text not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  ,@(\"2024-04-20 Saturday\", \"2024-04-21 Sunday\")\n  ,@(\"34\u00d785=2890\", \"92\u00d752=4784\")\n  ,@(\"84\u00d731=2604\", \"11\u00d746=506\")\n  ,@(\"91\u00d730=2730\", \"64\u00d769=4416\")\n  ,@(\"98\u00d779=7742\", \"47\u00d726=1222\")\n  ,@(\"59\u00d780=4720\", \"54\u00d787=4698\")\n  ,@(\"34\u00d740=1360\", \"15\u00d714=210\")\n  ,@(\"92\u00d716=1472\", \"48\u00d712=576\")\n  ,@(\"88\u00d795=8360\", \"56\u00d783=4648\")\n  ,@(\"60\u00d783=4980\", \"66\u00d719=1254\")\n  ,@(\"67\u00d784=5628\", \"90\u00d713=1170\")\n  ,@(\"64\u00d752=3328\", \"92\u00d725=2300\")\n  ,@(\"53\u00d729=1537\", \"95\u00d789=8455\")\n  ,@(\"38\u00d751=1938\", \"19\u00d735=665\")\n  ,@(\"26\u00d794=2444\", \"25\u00d784=2100\")\n  ,@(\"59\u00d796=5664\", \"82\u00d723=1886\")\n  ,@(\"42\u00d776=3192\", \"46\u00d749=2254\")\n  ,@(\"45\u00d725=1125\", \"84\u00d777=6468\")\n  ,@(\"38\u00d793=3534\", \"89\u00d733=2937\")\n  ,@(\"62\u00d785=5270\", \"60\u00d739=2340\")\n  ,@(\"47\u00d753=2491\", \"17\u00d756=952\")\n  ,@(\"51\u00d767=3417\", \"62\u00d790=5580\")\n  ,@(\"61\u00d749=2989\", \"95\u00d728=2660\")\n  ,@(\"29\u00d745=1305\", \"36\u00d737=1332\")\n  ,@(\"19\u00d746=874\", \"35\u00d798=3430\")\n  ,@(\"59\u00d787=5133\", \"27\u00d776=2052\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $find.Format = $false\n  $find.MatchCase = $true\n  $find.MatchWholeWord = $false\n  $find.MatchWildcards = $false\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}"}
